$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 58, shifting all rows
# from 58..77 down to 60..79 (data for those rows stays identical).
$ws.Rows("58:59").Insert()

# Row 58 - new weekly entry (Primera)
$ws.Cells.Item(58,1).Value  = 1
$ws.Cells.Item(58,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(58,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(58,4).Value  = 45135
$ws.Cells.Item(58,5).Value  = 15
$ws.Cells.Item(58,6).Value  = "Fruta"
$ws.Cells.Item(58,7).Value  = 100108
$ws.Cells.Item(58,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(58,9).Value  = 100108001
$ws.Cells.Item(58,10).Value = "Guayaba"
$ws.Cells.Item(58,11).Value = "Sin especificar"
$ws.Cells.Item(58,12).Value = "Primera"
$ws.Cells.Item(58,13).Value = 160
$ws.Cells.Item(58,14).Value = 4000
$ws.Cells.Item(58,15).Value = 4500
$ws.Cells.Item(58,16).Value = 4188
$ws.Cells.Item(58,17).Value = "$/caja 10 kilos"
$ws.Cells.Item(58,18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(58,19).Value = 419
$ws.Cells.Item(58,20).Value = 10

# Row 59 - new weekly entry (Segunda)
$ws.Cells.Item(59,1).Value  = 1
$ws.Cells.Item(59,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(59,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(59,4).Value  = 45135
$ws.Cells.Item(59,5).Value  = 15
$ws.Cells.Item(59,6).Value  = "Fruta"
$ws.Cells.Item(59,7).Value  = 100108
$ws.Cells.Item(59,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(59,9).Value  = 100108001
$ws.Cells.Item(59,10).Value = "Guayaba"
$ws.Cells.Item(59,11).Value = "Sin especificar"
$ws.Cells.Item(59,12).Value = "Segunda"
$ws.Cells.Item(59,13).Value = 150
$ws.Cells.Item(59,14).Value = 3500
$ws.Cells.Item(59,15).Value = 4000
$ws.Cells.Item(59,16).Value = 3833
$ws.Cells.Item(59,17).Value = "$/caja 10 kilos"
$ws.Cells.Item(59,18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(59,19).Value = 383
$ws.Cells.Item(59,20).Value = 10
